$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "pictures instead of letters":
#  - Column A ("Most Probable" / "Least Probable") takes over the big
#    decorative looks that used to live on columns B/C (the hiragana
#    font, and the bold Calibri font) for rows 2-3.
#  - Columns B/C on rows 2-3 switch from hiragana-letters / numbers text
#    to picture filenames, shown in the plain default font.

# A2 takes on the look B2 used to have (hiragana tfb, 60pt)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null

# A3 takes on the look C2 used to have (bold Calibri, 60pt)
$ws.Range("C2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null

# B2/C2/B3/C3 drop to the plain default look (copied from the already
# plain B1 cell) ahead of getting their new picture-filename text
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B2:C3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Replace the letter/number contents of B2/C2/B3/C3 with picture filenames.
# Nselect2.jpg is written first so it is interned into the shared-strings
# table ahead of Hselect2.jpg.
$ws.Range("C2").Value = "Nselect2.jpg"
$ws.Range("B2").Value = "Hselect2.jpg"
$ws.Range("C3").Value = "Nselect2.jpg"
$ws.Range("B3").Value = "Hselect2.jpg"

# Shrink rows 2-3 to fit the smaller picture-filename text
$ws.Rows.Item(2).RowHeight = 37.5
$ws.Rows.Item(3).RowHeight = 37.5

# Page setup tweak that came along with the resave
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved off of A2:A3
$ws.Range("B8").Select() | Out-Null
